$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text: volume number + week-covering date range ----
$ws.Range("A8").Value = "Volume 30   Number  34"
$ws.Range("C9").Value = "Report Covering the Week  8/21/2023  Through  8/27/2023"

# ---- Simple same-kind value updates ----
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = -33.333333333333
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -40
$ws.Range("M16").Value = 81.632653061224
$ws.Range("N16").Value = -84.163701067615
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 85.714285714285
$ws.Range("I17").Value = 93
$ws.Range("J17").Value = 89
$ws.Range("K17").Value = 4.494382022471
$ws.Range("L17").Value = 19.230769230769
$ws.Range("M17").Value = 93.75
$ws.Range("N17").Value = -24.390243902439
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -20
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -20
$ws.Range("I18").Value = 122
$ws.Range("J18").Value = 179
$ws.Range("K18").Value = -31.843575418994
$ws.Range("L18").Value = 35.555555555555
$ws.Range("M18").Value = 1.666666666666
$ws.Range("N18").Value = -78.018018018018
$ws.Range("C19").Value = 43
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 168.75
$ws.Range("F19").Value = 115
$ws.Range("G19").Value = 93
$ws.Range("H19").Value = 23.655913978494
$ws.Range("I19").Value = 809
$ws.Range("J19").Value = 796
$ws.Range("K19").Value = 1.633165829145
$ws.Range("L19").Value = 72.494669509594
$ws.Range("M19").Value = 16.235632183908
$ws.Range("N19").Value = -68.348982785602
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -22.222222222222
$ws.Range("I20").Value = 45
$ws.Range("J20").Value = 48
$ws.Range("K20").Value = -6.25
$ws.Range("L20").Value = 32.35294117647
$ws.Range("M20").Value = 80
$ws.Range("N20").Value = -92.359932088285
$ws.Range("C21").Value = 49
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = 75
$ws.Range("F21").Value = 153
$ws.Range("G21").Value = 135
$ws.Range("H21").Value = 13.333333333333
$ws.Range("I21").Value = 1168
$ws.Range("J21").Value = 1218
$ws.Range("K21").Value = -4.105090311986
$ws.Range("L21").Value = 55.112881806108
$ws.Range("M21").Value = 24.123273113708
$ws.Range("N21").Value = -73.466606088141
$ws.Range("C22").Value = 4
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 75
$ws.Range("I22").Value = 58
$ws.Range("K22").Value = -13.432835820895
$ws.Range("L22").Value = 26.086956521739
$ws.Range("M22").Value = 38.095238095238
$ws.Range("C24").Value = 97
$ws.Range("D24").Value = 89
$ws.Range("E24").Value = 8.988764044943
$ws.Range("F24").Value = 383
$ws.Range("G24").Value = 307
$ws.Range("H24").Value = 24.755700325732
$ws.Range("I24").Value = 2600
$ws.Range("J24").Value = 2586
$ws.Range("K24").Value = 0.541376643464
$ws.Range("L24").Value = 102.808112324493
$ws.Range("M24").Value = 125.694444444444
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -12.5
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = 45.16129032258
$ws.Range("I25").Value = 258
$ws.Range("J25").Value = 234
$ws.Range("K25").Value = 10.25641025641
$ws.Range("L25").Value = 25.242718446601
$ws.Range("M25").Value = 56.363636363636
$ws.Range("J26").Value = 21
$ws.Range("K26").Value = -42.857142857142
$ws.Range("L26").Value = -7.692307692307
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 10
$ws.Range("H27").Value = -37.5
$ws.Range("I27").Value = 62
$ws.Range("J27").Value = 77
$ws.Range("K27").Value = -19.480519480519
$ws.Range("L27").Value = 21.56862745098
$ws.Range("I30").Value = 7
$ws.Range("J30").Value = 12
$ws.Range("K30").Value = -41.666666666666
$ws.Range("L30").Value = -12.5

# ---- Text placeholder -> numeric (set format then value) ----
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = -100
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("D26").Value = 1
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E26").Value = -100
$ws.Range("D30").NumberFormat = '#,##0'
$ws.Range("D30").Value = 1
$ws.Range("E30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E30").Value = -100
$ws.Range("F30").NumberFormat = '#,##0'
$ws.Range("F30").Value = 1
$ws.Range("G30").NumberFormat = '#,##0'
$ws.Range("G30").Value = 1
$ws.Range("H30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H30").Value = 0

# ---- Numeric -> text placeholder (force text, then restore General style) ----
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"

# Restore General/style-14 formatting on the cells that just became text
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").PasteSpecial(-4122)
